$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 127.875
$ws.Range("I5").Value = 117.57143
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 117.57143
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -2.571430000000007
$ws.Range("N5").Value = -430

$ws.Range("H6").Value = 45.666668
$ws.Range("I6").Value = 45.666668
$ws.Range("K6").Value = 137.000004
$ws.Range("M6").Value = -25.00000399999999

$ws.Range("H9").Value = 317.14285
$ws.Range("J9").Value = 399.66666
$ws.Range("L9").Value = 399.66666
$ws.Range("N9").Value = -737.66666

$ws.Range("H12").Value = 831.8889
$ws.Range("I12").Value = 806.5
$ws.Range("K12").Value = 806.5
$ws.Range("M12").Value = -636.5

$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H29").Value = 4816.6665
$ws.Range("J29").Value = 6142.857
$ws.Range("L29").Value = 18428.571
$ws.Range("N29").Value = -18990.571

$ws.Range("H38").Value = 1442.3
$ws.Range("I38").Value = 47
$ws.Range("J38").Value = 14000
$ws.Range("K38").Value = 141
$ws.Range("L38").Value = 42000
$ws.Range("M38").Value = 231
$ws.Range("N38").Value = -42744

$ws.Range("H43").Value = 4458.4
$ws.Range("J43").Value = 4750
$ws.Range("L43").Value = 4750
$ws.Range("N43").Value = -4888

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H53").Value = 755.1818
$ws.Range("I53").Value = 460.69232
$ws.Range("J53").Value = 1180.5555
$ws.Range("K53").Value = 460.69232
$ws.Range("L53").Value = 1180.5555
$ws.Range("M53").Value = 176.30768
$ws.Range("N53").Value = -2454.5555

$ws.Range("H58").Value = 2388
$ws.Range("I58").Value = 864.8
$ws.Range("J58").Value = 3657.3333
$ws.Range("K58").Value = 2594.4
$ws.Range("L58").Value = 10971.9999
$ws.Range("M58").Value = -2444.4
$ws.Range("N58").Value = -11271.9999

$ws.Range("H132").Value = 30801.223
$ws.Range("I132").Value = 35368.5
$ws.Range("K132").Value = 106105.5
$ws.Range("M132").Value = -103575.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 499.75
$ws.Range("I23").Value = 499.75
$ws.Range("K23").Value = 499.75
$ws.Range("M23").Value = -240.75

$ws.Range("H24").Value = 2517500
$ws.Range("J24").Value = 2517500
$ws.Range("L24").Value = 2517500
$ws.Range("N24").Value = -2518248

$ws.Range("H44").Value = 12137.333
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 12137.333
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 12137.333
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -13113.333

$ws.Range("H55").Value = 45666.332
$ws.Range("J55").Value = 60999.5
$ws.Range("L55").Value = 60999.5
$ws.Range("N55").Value = -61629.5

$ws.Range("H100").Value = 2517500
$ws.Range("J100").Value = 2517500
$ws.Range("L100").Value = 2517500
$ws.Range("N100").Value = -2519664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3306
$ws.Range("I20").Value = 3478.8
$ws.Range("J20").Value = 3018
$ws.Range("K20").Value = 3478.8
$ws.Range("L20").Value = 3018
$ws.Range("M20").Value = -3231.8
$ws.Range("N20").Value = -3512

$ws.Range("H62").Value = 55714.145
$ws.Range("J62").Value = 48333.332
$ws.Range("L62").Value = 48333.332
$ws.Range("N62").Value = -49705.332

$ws.Range("H65").Value = 55714.145
$ws.Range("J65").Value = 48333.332
$ws.Range("L65").Value = 144999.996
$ws.Range("N65").Value = -151863.996

$ws.Range("H86").Value = 4965.7856
$ws.Range("I86").Value = 1739
$ws.Range("J86").Value = 7385.875
$ws.Range("K86").Value = 1739
$ws.Range("L86").Value = 7385.875
$ws.Range("M86").Value = -616
$ws.Range("N86").Value = -9631.875

$ws.Range("H89").Value = 4965.7856
$ws.Range("I89").Value = 1739
$ws.Range("J89").Value = 7385.875
$ws.Range("K89").Value = 8695
$ws.Range("L89").Value = 36929.375
$ws.Range("M89").Value = -3079
$ws.Range("N89").Value = -48161.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 48571.285
$ws.Range("J63").Value = 48571.285
$ws.Range("L63").Value = 48571.285
$ws.Range("N63").Value = -49943.285

$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

$ws.Range("H66").Value = 48571.285
$ws.Range("J66").Value = 48571.285
$ws.Range("L66").Value = 145713.855
$ws.Range("N66").Value = -152577.855

$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws.Range("H99").Value = 2480
$ws.Range("I99").Value = 2478.25
$ws.Range("J99").Value = 2481.75
$ws.Range("K99").Value = 2478.25
$ws.Range("L99").Value = 2481.75
$ws.Range("M99").Value = -980.25
$ws.Range("N99").Value = -5477.75

$ws.Range("H126").Value = 2480
$ws.Range("I126").Value = 2478.25
$ws.Range("J126").Value = 2481.75
$ws.Range("K126").Value = 7434.75
$ws.Range("L126").Value = 7445.25
$ws.Range("M126").Value = -4964.75
$ws.Range("N126").Value = -12385.25

$ws.Range("H134").Value = 3611.3333
$ws.Range("I134").Value = 3611.3333
$ws.Range("K134").Value = 10833.9999
$ws.Range("M134").Value = -8298.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2250
$ws.Range("I134").Value = 2250
$ws.Range("K134").Value = 6750
$ws.Range("M134").Value = -1680

$ws.Range("H137").Value = 3175
$ws.Range("I137").Value = 1875
$ws.Range("K137").Value = 5625
$ws.Range("M137").Value = -525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 321.73685
$ws.Range("I2").Value = 148.7
$ws.Range("K2").Value = 148.7
$ws.Range("M2").Value = -35.69999999999999

$ws.Range("H14").Value = 600
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1336

$ws.Range("H23").Value = 733.3333
$ws.Range("J23").Value = 733.3333
$ws.Range("L23").Value = 733.3333
$ws.Range("N23").Value = -1179.3333

$ws.Range("H42").Value = 98000
$ws.Range("J42").Value = 98000
$ws.Range("L42").Value = 98000
$ws.Range("N42").Value = -98970

$ws.Range("H43").Value = 19142.715
$ws.Range("J43").Value = 19833.166
$ws.Range("L43").Value = 19833.166
$ws.Range("N43").Value = -20135.166

$ws.Range("H46").Value = 3000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H57").Value = 61610
$ws.Range("J57").Value = 61610
$ws.Range("L57").Value = 61610
$ws.Range("N57").Value = -63250

$ws.Range("H113").Value = 3282.6667
$ws.Range("J113").Value = 4199
$ws.Range("L113").Value = 4199
$ws.Range("N113").Value = -8539

$ws.Range("H115").Value = 98000
$ws.Range("J115").Value = 98000
$ws.Range("L115").Value = 98000
$ws.Range("N115").Value = -100350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2375.5
$ws.Range("J16").Value = 2001
$ws.Range("L16").Value = 2001
$ws.Range("N16").Value = -2341

$ws.Range("H22").Value = 1216.5
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705

$ws.Range("H27").Value = 1216.5
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20450

$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21560

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 2700
$ws.Range("J20").Value = 2700
$ws.Range("L20").Value = 2700
$ws.Range("N20").Value = -3180

$ws.Range("H63").Value = 10000
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11248

$ws.Range("H66").Value = 10000
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36240

$ws.Range("I132").Value = 1737.0454
$ws.Range("J132").Value = 2066.6667
$ws.Range("K132").Value = 5211.1362
$ws.Range("L132").Value = 6200.000100000001
$ws.Range("M132").Value = -2681.1362
$ws.Range("N132").Value = -11260.0001

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
